$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44791
$ws.Cells.Item(2, 11).Value = 20000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 20000
$ws.Cells.Item(2, 16).Value = 1333
$ws.Cells.Item(3, 4).Value = 44446
$ws.Cells.Item(3, 10).Value = 34
$ws.Cells.Item(3, 11).Value = 24000
$ws.Cells.Item(3, 12).Value = 25000
$ws.Cells.Item(3, 13).Value = 24500
$ws.Cells.Item(3, 16).Value = 1633
$ws.Cells.Item(4, 4).Value = 44411
$ws.Cells.Item(4, 10).Value = 34
$ws.Cells.Item(4, 11).Value = 25000
$ws.Cells.Item(4, 12).Value = 26000
$ws.Cells.Item(4, 13).Value = 25500
$ws.Cells.Item(4, 16).Value = 1700
$ws.Cells.Item(5, 4).Value = 44343
$ws.Cells.Item(5, 10).Value = 26
$ws.Cells.Item(5, 11).Value = 23000
$ws.Cells.Item(5, 12).Value = 24000
$ws.Cells.Item(5, 13).Value = 23500
$ws.Cells.Item(5, 16).Value = 1567
$ws.Cells.Item(6, 4).Value = 44784
$ws.Cells.Item(6, 10).Value = 28
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 21000
$ws.Cells.Item(6, 13).Value = 20357
$ws.Cells.Item(6, 16).Value = 1357
$ws.Cells.Item(7, 4).Value = 44778
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 18000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 18000
$ws.Cells.Item(7, 16).Value = 1200
$ws.Cells.Item(8, 4).Value = 44453
$ws.Cells.Item(8, 10).Value = 25
$ws.Cells.Item(8, 11).Value = 25000
$ws.Cells.Item(8, 12).Value = 26000
$ws.Cells.Item(8, 13).Value = 25520
$ws.Cells.Item(8, 15).Value = "Hijuelas"
$ws.Cells.Item(8, 16).Value = 1701
$ws.Cells.Item(9, 4).Value = 45084
$ws.Cells.Item(9, 10).Value = 43
$ws.Cells.Item(9, 11).Value = 19000
$ws.Cells.Item(9, 12).Value = 21000
$ws.Cells.Item(9, 13).Value = 20023
$ws.Cells.Item(9, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value = 1335
$ws.Cells.Item(10, 4).Value = 44425
$ws.Cells.Item(10, 10).Value = 25
$ws.Cells.Item(10, 11).Value = 24000
$ws.Cells.Item(10, 12).Value = 25000
$ws.Cells.Item(10, 13).Value = 24520
$ws.Cells.Item(10, 16).Value = 1635
$ws.Cells.Item(11, 4).Value = 44349
$ws.Cells.Item(11, 10).Value = 21
$ws.Cells.Item(11, 11).Value = 24000
$ws.Cells.Item(11, 12).Value = 25000
$ws.Cells.Item(11, 13).Value = 24524
$ws.Cells.Item(11, 16).Value = 1635
$ws.Cells.Item(12, 4).Value = 44804
$ws.Cells.Item(12, 10).Value = 35
$ws.Cells.Item(12, 11).Value = 19000
$ws.Cells.Item(12, 12).Value = 20000
$ws.Cells.Item(12, 13).Value = 19714
$ws.Cells.Item(12, 16).Value = 1314
$ws.Cells.Item(13, 4).Value = 44771
$ws.Cells.Item(13, 10).Value = 43
$ws.Cells.Item(13, 11).Value = 22000
$ws.Cells.Item(13, 12).Value = 22000
$ws.Cells.Item(13, 13).Value = 22000
$ws.Cells.Item(13, 16).Value = 1467
$ws.Cells.Item(14, 4).Value = 45086
$ws.Cells.Item(14, 10).Value = 70
$ws.Cells.Item(14, 11).Value = 16000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 17000
$ws.Cells.Item(14, 16).Value = 1133
$ws.Cells.Item(15, 4).Value = 45100
$ws.Cells.Item(15, 10).Value = 43
$ws.Cells.Item(15, 11).Value = 17000
$ws.Cells.Item(15, 12).Value = 18000
$ws.Cells.Item(15, 13).Value = 17512
$ws.Cells.Item(15, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 16).Value = 1167
$ws.Cells.Item(16, 4).Value = 44838
$ws.Cells.Item(16, 10).Value = 52
$ws.Cells.Item(16, 11).Value = 17000
$ws.Cells.Item(16, 12).Value = 17000
$ws.Cells.Item(16, 13).Value = 17000
$ws.Cells.Item(16, 16).Value = 1133
$ws.Cells.Item(17, 4).Value = 44719
$ws.Cells.Item(17, 10).Value = 43
$ws.Cells.Item(17, 11).Value = 17000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 17512
$ws.Cells.Item(17, 16).Value = 1167
$ws.Cells.Item(18, 4).Value = 44749
$ws.Cells.Item(18, 10).Value = 34
$ws.Cells.Item(18, 11).Value = 18000
$ws.Cells.Item(18, 12).Value = 18000
$ws.Cells.Item(18, 13).Value = 18000
$ws.Cells.Item(18, 16).Value = 1200
$ws.Cells.Item(19, 4).Value = 44727
$ws.Cells.Item(19, 10).Value = 28
$ws.Cells.Item(19, 11).Value = 24000
$ws.Cells.Item(19, 13).Value = 24000
$ws.Cells.Item(19, 15).Value = "Hijuelas"
$ws.Cells.Item(19, 16).Value = 1600
$ws.Cells.Item(20, 4).Value = 44811
$ws.Cells.Item(20, 10).Value = 18
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 15).Value = "Hijuelas"
$ws.Cells.Item(20, 16).Value = 1333
$ws.Cells.Item(21, 4).Value = 45044
$ws.Cells.Item(21, 10).Value = 52
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 22000
$ws.Cells.Item(21, 13).Value = 21000
$ws.Cells.Item(21, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(21, 16).Value = 1400
$ws.Cells.Item(22, 4).Value = 44329
$ws.Cells.Item(22, 10).Value = 25
$ws.Cells.Item(22, 11).Value = 23000
$ws.Cells.Item(22, 12).Value = 23000
$ws.Cells.Item(22, 13).Value = 23000
$ws.Cells.Item(22, 15).Value = "Hijuelas"
$ws.Cells.Item(22, 16).Value = 1533
$ws.Cells.Item(23, 4).Value = 44757
$ws.Cells.Item(23, 10).Value = 34
$ws.Cells.Item(23, 11).Value = 17000
$ws.Cells.Item(23, 12).Value = 18000
$ws.Cells.Item(23, 13).Value = 17500
$ws.Cells.Item(23, 16).Value = 1167
$ws.Cells.Item(24, 4).Value = 44754
$ws.Cells.Item(24, 10).Value = 43
$ws.Cells.Item(24, 11).Value = 22000
$ws.Cells.Item(24, 12).Value = 22000
$ws.Cells.Item(24, 13).Value = 22000
$ws.Cells.Item(24, 16).Value = 1467
$ws.Cells.Item(25, 4).Value = 44460
$ws.Cells.Item(26, 4).Value = 45070
$ws.Cells.Item(26, 10).Value = 70
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 22000
$ws.Cells.Item(26, 13).Value = 21000
$ws.Cells.Item(26, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(26, 16).Value = 1400
$ws.Cells.Item(27, 4).Value = 44677
$ws.Cells.Item(27, 10).Value = 34
$ws.Cells.Item(27, 11).Value = 25000
$ws.Cells.Item(27, 12).Value = 26000
$ws.Cells.Item(27, 13).Value = 25500
$ws.Cells.Item(27, 16).Value = 1700
$ws.Cells.Item(28, 4).Value = 44783
$ws.Cells.Item(28, 10).Value = 35
$ws.Cells.Item(28, 11).Value = 20000
$ws.Cells.Item(28, 12).Value = 21000
$ws.Cells.Item(28, 13).Value = 20429
$ws.Cells.Item(28, 16).Value = 1362
$ws.Cells.Item(29, 4).Value = 44406
$ws.Cells.Item(29, 10).Value = 25
$ws.Cells.Item(29, 11).Value = 24000
$ws.Cells.Item(29, 12).Value = 25000
$ws.Cells.Item(29, 13).Value = 24520
$ws.Cells.Item(29, 16).Value = 1635
$ws.Cells.Item(30, 4).Value = 45063
$ws.Cells.Item(30, 10).Value = 52
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 22000
$ws.Cells.Item(30, 13).Value = 21000
$ws.Cells.Item(30, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(30, 16).Value = 1400
$ws.Cells.Item(31, 4).Value = 45022
$ws.Cells.Item(31, 10).Value = 15
$ws.Cells.Item(31, 11).Value = 27000
$ws.Cells.Item(31, 12).Value = 27000
$ws.Cells.Item(31, 13).Value = 27000
$ws.Cells.Item(31, 14).Value = "`$/malla 17 kilos"
$ws.Cells.Item(31, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(31, 16).Value = 1588
$ws.Cells.Item(31, 17).Value = 17
$ws.Cells.Item(32, 4).Value = 44819
$ws.Cells.Item(32, 11).Value = 14000
$ws.Cells.Item(32, 12).Value = 15000
$ws.Cells.Item(32, 13).Value = 14500
$ws.Cells.Item(32, 15).Value = "Hijuelas"
$ws.Cells.Item(32, 16).Value = 967
$ws.Cells.Item(33, 4).Value = 45091
$ws.Cells.Item(33, 10).Value = 34
$ws.Cells.Item(33, 11).Value = 19000
$ws.Cells.Item(33, 13).Value = 20000
$ws.Cells.Item(33, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(33, 16).Value = 1333
$ws.Cells.Item(34, 4).Value = 44385
$ws.Cells.Item(34, 10).Value = 25
$ws.Cells.Item(34, 11).Value = 14000
$ws.Cells.Item(34, 12).Value = 15000
$ws.Cells.Item(34, 13).Value = 14480
$ws.Cells.Item(34, 15).Value = "Hijuelas"
$ws.Cells.Item(34, 16).Value = 965
$ws.Cells.Item(35, 4).Value = 45055
$ws.Cells.Item(35, 10).Value = 52
$ws.Cells.Item(35, 11).Value = 22000
$ws.Cells.Item(35, 12).Value = 24000
$ws.Cells.Item(35, 13).Value = 23000
$ws.Cells.Item(35, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(35, 16).Value = 1533
$ws.Cells.Item(36, 4).Value = 44413
$ws.Cells.Item(36, 10).Value = 25
$ws.Cells.Item(36, 11).Value = 24000
$ws.Cells.Item(36, 12).Value = 25000
$ws.Cells.Item(36, 13).Value = 24480
$ws.Cells.Item(36, 15).Value = "Hijuelas"
$ws.Cells.Item(36, 16).Value = 1632
$ws.Cells.Item(38, 4).Value = 45072
$ws.Cells.Item(38, 10).Value = 43
$ws.Cells.Item(38, 11).Value = 22000
$ws.Cells.Item(38, 12).Value = 24000
$ws.Cells.Item(38, 13).Value = 23023
$ws.Cells.Item(38, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(38, 16).Value = 1535
$ws.Cells.Item(39, 4).Value = 44726
$ws.Cells.Item(39, 10).Value = 28
$ws.Cells.Item(39, 12).Value = 24000
$ws.Cells.Item(39, 13).Value = 24000
$ws.Cells.Item(39, 16).Value = 1600
$ws.Cells.Item(40, 4).Value = 44792
$ws.Cells.Item(40, 10).Value = 56
$ws.Cells.Item(40, 11).Value = 19000
$ws.Cells.Item(40, 13).Value = 19500
$ws.Cells.Item(40, 16).Value = 1300
$ws.Cells.Item(41, 4).Value = 45079
$ws.Cells.Item(41, 10).Value = 52
$ws.Cells.Item(41, 11).Value = 19000
$ws.Cells.Item(41, 12).Value = 21000
$ws.Cells.Item(41, 13).Value = 20000
$ws.Cells.Item(41, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(41, 16).Value = 1333
$ws.Cells.Item(42, 4).Value = 44750
$ws.Cells.Item(42, 12).Value = 22000
$ws.Cells.Item(42, 13).Value = 21000
$ws.Cells.Item(42, 16).Value = 1400
$ws.Cells.Item(43, 4).Value = 44831
$ws.Cells.Item(43, 10).Value = 20
$ws.Cells.Item(43, 12).Value = 20000
$ws.Cells.Item(43, 13).Value = 19500
$ws.Cells.Item(43, 15).Value = "Hijuelas"
$ws.Cells.Item(43, 16).Value = 1300
$ws.Cells.Item(44, 4).Value = 44817
$ws.Cells.Item(44, 11).Value = 20000
$ws.Cells.Item(44, 12).Value = 20000
$ws.Cells.Item(44, 13).Value = 20000
$ws.Cells.Item(44, 16).Value = 1333
$ws.Cells.Item(45, 4).Value = 44707
$ws.Cells.Item(45, 10).Value = 30
$ws.Cells.Item(45, 11).Value = 26000
$ws.Cells.Item(45, 13).Value = 26000
$ws.Cells.Item(45, 16).Value = 1733
$ws.Cells.Item(46, 4).Value = 45083
$ws.Cells.Item(46, 10).Value = 34
$ws.Cells.Item(46, 11).Value = 19000
$ws.Cells.Item(46, 12).Value = 21000
$ws.Cells.Item(46, 13).Value = 20000
$ws.Cells.Item(46, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(46, 16).Value = 1333
$ws.Cells.Item(47, 4).Value = 44400
$ws.Cells.Item(47, 10).Value = 16
$ws.Cells.Item(48, 4).Value = 44435
$ws.Cells.Item(48, 10).Value = 34
$ws.Cells.Item(48, 12).Value = 25000
$ws.Cells.Item(48, 13).Value = 24500
$ws.Cells.Item(48, 16).Value = 1633
$ws.Cells.Item(49, 4).Value = 44797
$ws.Cells.Item(49, 10).Value = 43
$ws.Cells.Item(49, 11).Value = 20000
$ws.Cells.Item(49, 12).Value = 20000
$ws.Cells.Item(49, 13).Value = 20000
$ws.Cells.Item(49, 15).Value = "Hijuelas"
$ws.Cells.Item(49, 16).Value = 1333
$ws.Cells.Item(50, 4).Value = 44336
$ws.Cells.Item(50, 10).Value = 34
$ws.Cells.Item(50, 13).Value = 24500
$ws.Cells.Item(50, 16).Value = 1633
$ws.Cells.Item(51, 4).Value = 44761
$ws.Cells.Item(51, 10).Value = 43
$ws.Cells.Item(51, 11).Value = 19000
$ws.Cells.Item(51, 12).Value = 19000
$ws.Cells.Item(51, 13).Value = 19000
$ws.Cells.Item(51, 16).Value = 1267
$ws.Cells.Item(52, 4).Value = 45085
$ws.Cells.Item(52, 10).Value = 52
$ws.Cells.Item(52, 11).Value = 20000
$ws.Cells.Item(52, 12).Value = 22000
$ws.Cells.Item(52, 13).Value = 21000
$ws.Cells.Item(52, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(52, 16).Value = 1400
$ws.Cells.Item(53, 4).Value = 44736
$ws.Cells.Item(53, 10).Value = 27
$ws.Cells.Item(53, 11).Value = 24000
$ws.Cells.Item(53, 12).Value = 25000
$ws.Cells.Item(53, 13).Value = 24481
$ws.Cells.Item(53, 16).Value = 1632
$ws.Cells.Item(54, 4).Value = 44841
$ws.Cells.Item(54, 10).Value = 45
$ws.Cells.Item(54, 11).Value = 17000
$ws.Cells.Item(54, 12).Value = 17000
$ws.Cells.Item(54, 13).Value = 17000
$ws.Cells.Item(54, 16).Value = 1133
$ws.Cells.Item(55, 4).Value = 44390
$ws.Cells.Item(55, 10).Value = 34
$ws.Cells.Item(55, 11).Value = 24000
$ws.Cells.Item(55, 12).Value = 25000
$ws.Cells.Item(55, 13).Value = 24500
$ws.Cells.Item(55, 16).Value = 1633
$ws.Cells.Item(56, 4).Value = 44432
$ws.Cells.Item(56, 10).Value = 34
$ws.Cells.Item(56, 11).Value = 24000
$ws.Cells.Item(56, 12).Value = 25000
$ws.Cells.Item(56, 13).Value = 24500
$ws.Cells.Item(56, 16).Value = 1633
$ws.Cells.Item(57, 4).Value = 44818
$ws.Cells.Item(57, 11).Value = 20000
$ws.Cells.Item(57, 12).Value = 20000
$ws.Cells.Item(57, 15).Value = "Hijuelas"
$ws.Cells.Item(58, 4).Value = 44418
$ws.Cells.Item(58, 10).Value = 16
$ws.Cells.Item(58, 11).Value = 25000
$ws.Cells.Item(58, 12).Value = 26000
$ws.Cells.Item(58, 13).Value = 25500
$ws.Cells.Item(58, 16).Value = 1700
$ws.Cells.Item(59, 4).Value = 44463
$ws.Cells.Item(59, 10).Value = 25
$ws.Cells.Item(59, 11).Value = 24000
$ws.Cells.Item(59, 12).Value = 25000
$ws.Cells.Item(59, 13).Value = 24480
$ws.Cells.Item(59, 16).Value = 1632
$ws.Cells.Item(60, 4).Value = 44806
$ws.Cells.Item(60, 10).Value = 27
$ws.Cells.Item(60, 11).Value = 19000
$ws.Cells.Item(60, 12).Value = 20000
$ws.Cells.Item(60, 13).Value = 19556
$ws.Cells.Item(60, 16).Value = 1304
$ws.Cells.Item(61, 4).Value = 45092
$ws.Cells.Item(61, 10).Value = 43
$ws.Cells.Item(61, 11).Value = 17000
$ws.Cells.Item(61, 12).Value = 19000
$ws.Cells.Item(61, 13).Value = 18023
$ws.Cells.Item(61, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(61, 16).Value = 1202
$ws.Cells.Item(61, 17).Value = 15
$ws.Cells.Item(62, 4).Value = 44832
$ws.Cells.Item(62, 10).Value = 22
$ws.Cells.Item(62, 11).Value = 20000
$ws.Cells.Item(62, 12).Value = 20000
$ws.Cells.Item(62, 13).Value = 20000
$ws.Cells.Item(62, 16).Value = 1333
$ws.Cells.Item(63, 4).Value = 45062
$ws.Cells.Item(63, 11).Value = 22000
$ws.Cells.Item(63, 12).Value = 25000
$ws.Cells.Item(63, 13).Value = 23500
$ws.Cells.Item(63, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(63, 16).Value = 1567
$ws.Cells.Item(64, 4).Value = 44708
$ws.Cells.Item(64, 10).Value = 25
$ws.Cells.Item(64, 11).Value = 26000
$ws.Cells.Item(64, 12).Value = 26000
$ws.Cells.Item(64, 13).Value = 26000
$ws.Cells.Item(64, 16).Value = 1733
$ws.Cells.Item(65, 4).Value = 44351
$ws.Cells.Item(65, 10).Value = 34
$ws.Cells.Item(65, 12).Value = 25000
$ws.Cells.Item(65, 13).Value = 24500
$ws.Cells.Item(65, 16).Value = 1633
$ws.Cells.Item(66, 4).Value = 45050
$ws.Cells.Item(66, 10).Value = 70
$ws.Cells.Item(66, 11).Value = 20000
$ws.Cells.Item(66, 12).Value = 22000
$ws.Cells.Item(66, 13).Value = 21000
$ws.Cells.Item(66, 16).Value = 1400
$ws.Cells.Item(67, 4).Value = 45071
$ws.Cells.Item(67, 10).Value = 52
$ws.Cells.Item(67, 12).Value = 24000
$ws.Cells.Item(67, 13).Value = 23000
$ws.Cells.Item(67, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(67, 16).Value = 1533
$ws.Cells.Item(68, 4).Value = 44789
$ws.Cells.Item(68, 10).Value = 34
$ws.Cells.Item(68, 11).Value = 21000
$ws.Cells.Item(68, 12).Value = 21000
$ws.Cells.Item(68, 13).Value = 21000
$ws.Cells.Item(68, 16).Value = 1400
$ws.Cells.Item(69, 4).Value = 44685
$ws.Cells.Item(69, 10).Value = 20
$ws.Cells.Item(69, 11).Value = 25000
$ws.Cells.Item(69, 13).Value = 25000
$ws.Cells.Item(69, 15).Value = "Hijuelas"
$ws.Cells.Item(69, 16).Value = 1667
$ws.Cells.Item(70, 4).Value = 44729
$ws.Cells.Item(70, 10).Value = 52
$ws.Cells.Item(70, 11).Value = 24000
$ws.Cells.Item(70, 12).Value = 24000
$ws.Cells.Item(70, 13).Value = 24000
$ws.Cells.Item(70, 16).Value = 1600
$ws.Cells.Item(71, 4).Value = 44421
$ws.Cells.Item(71, 10).Value = 18
$ws.Cells.Item(71, 11).Value = 24000
$ws.Cells.Item(71, 12).Value = 25000
$ws.Cells.Item(71, 13).Value = 24500
$ws.Cells.Item(71, 16).Value = 1633
$ws.Cells.Item(72, 4).Value = 44776
$ws.Cells.Item(72, 10).Value = 34
$ws.Cells.Item(72, 11).Value = 20000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 20000
$ws.Cells.Item(72, 16).Value = 1333
$ws.Cells.Item(73, 4).Value = 44747
$ws.Cells.Item(74, 4).Value = 44790
$ws.Cells.Item(74, 10).Value = 36
$ws.Cells.Item(74, 12).Value = 20000
$ws.Cells.Item(74, 13).Value = 20000
$ws.Cells.Item(74, 16).Value = 1333
$ws.Cells.Item(75, 4).Value = 44775
$ws.Cells.Item(75, 10).Value = 43
$ws.Cells.Item(75, 11).Value = 20000
$ws.Cells.Item(75, 12).Value = 20000
$ws.Cells.Item(75, 13).Value = 20000
$ws.Cells.Item(75, 15).Value = "Hijuelas"
$ws.Cells.Item(75, 16).Value = 1333
$ws.Cells.Item(76, 4).Value = 45093
$ws.Cells.Item(76, 10).Value = 52
$ws.Cells.Item(76, 11).Value = 17000
$ws.Cells.Item(76, 12).Value = 19000
$ws.Cells.Item(76, 13).Value = 18000
$ws.Cells.Item(76, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(76, 16).Value = 1200
$ws.Cells.Item(77, 4).Value = 44714
$ws.Cells.Item(77, 10).Value = 52
$ws.Cells.Item(77, 11).Value = 18000
$ws.Cells.Item(77, 12).Value = 20000
$ws.Cells.Item(77, 13).Value = 19000
$ws.Cells.Item(77, 16).Value = 1267
$ws.Cells.Item(78, 4).Value = 44763
$ws.Cells.Item(78, 10).Value = 34
$ws.Cells.Item(78, 12).Value = 20000
$ws.Cells.Item(78, 13).Value = 20000
$ws.Cells.Item(78, 15).Value = "Hijuelas"
$ws.Cells.Item(78, 16).Value = 1333
$ws.Cells.Item(79, 4).Value = 44706
$ws.Cells.Item(79, 10).Value = 30
$ws.Cells.Item(79, 11).Value = 26000
$ws.Cells.Item(79, 12).Value = 26000
$ws.Cells.Item(79, 13).Value = 26000
$ws.Cells.Item(79, 16).Value = 1733
$ws.Cells.Item(80, 4).Value = 44742
$ws.Cells.Item(80, 11).Value = 20000
$ws.Cells.Item(80, 12).Value = 22000
$ws.Cells.Item(80, 13).Value = 21000
$ws.Cells.Item(80, 16).Value = 1400
$ws.Cells.Item(81, 4).Value = 44705
$ws.Cells.Item(81, 10).Value = 35
$ws.Cells.Item(81, 11).Value = 26000
$ws.Cells.Item(81, 12).Value = 26000
$ws.Cells.Item(81, 13).Value = 26000
$ws.Cells.Item(81, 16).Value = 1733
$ws.Cells.Item(82, 4).Value = 44680
$ws.Cells.Item(82, 10).Value = 36
$ws.Cells.Item(82, 11).Value = 24000
$ws.Cells.Item(82, 12).Value = 25000
$ws.Cells.Item(82, 13).Value = 24500
$ws.Cells.Item(82, 16).Value = 1633
$ws.Cells.Item(83, 4).Value = 44397
$ws.Cells.Item(83, 10).Value = 34
$ws.Cells.Item(83, 11).Value = 23000
$ws.Cells.Item(83, 12).Value = 24000
$ws.Cells.Item(83, 13).Value = 23500
$ws.Cells.Item(83, 15).Value = "Hijuelas"
$ws.Cells.Item(83, 16).Value = 1567
$ws.Cells.Item(84, 4).Value = 44428
$ws.Cells.Item(84, 10).Value = 16
$ws.Cells.Item(84, 11).Value = 25000
$ws.Cells.Item(84, 12).Value = 26000
$ws.Cells.Item(84, 13).Value = 25500
$ws.Cells.Item(84, 16).Value = 1700
$ws.Cells.Item(85, 4).Value = 44799
$ws.Cells.Item(85, 10).Value = 27
$ws.Cells.Item(85, 11).Value = 20000
$ws.Cells.Item(85, 12).Value = 20000
$ws.Cells.Item(85, 13).Value = 20000
$ws.Cells.Item(85, 16).Value = 1333
$ws.Cells.Item(86, 4).Value = 44755
$ws.Cells.Item(86, 10).Value = 43
$ws.Cells.Item(86, 11).Value = 18000
$ws.Cells.Item(86, 12).Value = 20000
$ws.Cells.Item(86, 13).Value = 19023
$ws.Cells.Item(86, 16).Value = 1268
$ws.Cells.Item(87, 4).Value = 44341
$ws.Cells.Item(87, 10).Value = 36
$ws.Cells.Item(87, 11).Value = 24000
$ws.Cells.Item(87, 12).Value = 25000
$ws.Cells.Item(87, 13).Value = 24500
$ws.Cells.Item(87, 16).Value = 1633
$ws.Cells.Item(88, 4).Value = 44383
$ws.Cells.Item(88, 11).Value = 13000
$ws.Cells.Item(88, 12).Value = 14000
$ws.Cells.Item(88, 13).Value = 13480
$ws.Cells.Item(88, 16).Value = 899
$ws.Cells.Item(89, 4).Value = 44810
$ws.Cells.Item(89, 10).Value = 22
$ws.Cells.Item(89, 11).Value = 20000
$ws.Cells.Item(89, 13).Value = 20000
$ws.Cells.Item(89, 16).Value = 1333
$ws.Cells.Item(90, 4).Value = 44442
$ws.Cells.Item(90, 10).Value = 28
$ws.Cells.Item(90, 11).Value = 24000
$ws.Cells.Item(90, 12).Value = 25000
$ws.Cells.Item(90, 13).Value = 24500
$ws.Cells.Item(90, 15).Value = "Hijuelas"
$ws.Cells.Item(90, 16).Value = 1633
$ws.Cells.Item(91, 4).Value = 44769
$ws.Cells.Item(91, 10).Value = 34
$ws.Cells.Item(91, 11).Value = 20000
$ws.Cells.Item(91, 12).Value = 20000
$ws.Cells.Item(91, 13).Value = 20000
$ws.Cells.Item(91, 16).Value = 1333
$ws.Cells.Item(92, 4).Value = 44449
$ws.Cells.Item(92, 10).Value = 18
$ws.Cells.Item(93, 4).Value = 44455
$ws.Cells.Item(93, 10).Value = 18
$ws.Cells.Item(93, 11).Value = 24000
$ws.Cells.Item(93, 12).Value = 25000
$ws.Cells.Item(93, 13).Value = 24500
$ws.Cells.Item(93, 16).Value = 1633
